$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.198.50'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.852.74'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''0.7002'
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").Value = '''237.25'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.07884'
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("D9").Value = '''0.3018'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").Value = '''24.12'
$ws.Range("E10").Value = '  +4.40%  '
$ws.Range("D11").Value = '''0.08147'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '1.947.43'
$ws.Range("E12").Value = '  +4.94%  '
$ws.Range("D13").Value = '''5.185'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '''0.7049'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").Value = '''89.32'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '29.329.70'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '''5.778'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").Value = '''0.000007840'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = '''13.19'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '''235.63'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.140.93'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''7.510'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").Value = '''161.77'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '''8.881'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = '''0.1419'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").Value = '''18.03'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").Value = '''1.900'
$ws.Range("E29").Value = '  -2.45%  '
$ws.Range("D30").Value = '''1.398'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").Value = '''1.473'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = '''4.284'
$ws.Range("E32").Value = '  -4.80%  '
$ws.Range("D33").Value = '''4.011'
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '''0.05147'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").Value = '''1.167'
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").Value = '''0.7065'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").Value = '''1.000'
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").Value = '''2.677'
$ws.Range("D39").Value = '''0.01846'
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").Value = '''2.697'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '1.149.22'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("D42").Value = '''0.9202'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").Value = '''5.951'
$ws.Range("E43").Value = '  +0.90%  '
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '''102.73'
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '''0.5287'
$ws.Range("E48").Value = '  -3.72%  '
$ws.Range("D49").Value = '''1.735'
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.149'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '''6.950'
$ws.Range("E51").Value = '  -0.62%  '
